$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 852: simplify DragAttachTag messages into a single "old -> new" format
$ws.Range("B852").Value = "标签: %s -> %s"
$ws.Range("C852").Value = "Tag: %s -> %s"

# Rows 853-873: new PropertyBrush property rows (key / Chinese / English)
$ws.Range("A853").Value = "PropertyBrush.House"
$ws.Range("B853").Value = "所属方"
$ws.Range("C853").Value = "House"
$ws.Range("A854").Value = "PropertyBrush.Health"
$ws.Range("B854").Value = "强度"
$ws.Range("C854").Value = "Health"
$ws.Range("A855").Value = "PropertyBrush.Facing"
$ws.Range("B855").Value = "方向"
$ws.Range("C855").Value = "Facing"
$ws.Range("A856").Value = "PropertyBrush.AISellable"
$ws.Range("B856").Value = "AI变卖"
$ws.Range("C856").Value = "AISellable"
$ws.Range("A857").Value = "PropertyBrush.AIRebuildable"
$ws.Range("B857").Value = "重建"
$ws.Range("C857").Value = "AIRebuildable"
$ws.Range("A858").Value = "PropertyBrush.PoweredOn"
$ws.Range("B858").Value = "耗能/工作"
$ws.Range("C858").Value = "PoweredOn"
$ws.Range("A859").Value = "PropertyBrush.Upgrades"
$ws.Range("B859").Value = "组件数"
$ws.Range("C859").Value = "Upgrades"
$ws.Range("A860").Value = "PropertyBrush.SpotLight"
$ws.Range("B860").Value = "探照灯"
$ws.Range("C860").Value = "SpotLight"
$ws.Range("A861").Value = "PropertyBrush.Upgrade1"
$ws.Range("B861").Value = "组件1"
$ws.Range("C861").Value = "Upgrade1"
$ws.Range("A862").Value = "PropertyBrush.Upgrade2"
$ws.Range("B862").Value = "组件2"
$ws.Range("C862").Value = "Upgrade2"
$ws.Range("A863").Value = "PropertyBrush.Upgrade3"
$ws.Range("B863").Value = "组件3"
$ws.Range("C863").Value = "Upgrade3"
$ws.Range("A864").Value = "PropertyBrush.AIRepairable"
$ws.Range("B864").Value = "AI修复"
$ws.Range("C864").Value = "AIRepairable"
$ws.Range("A865").Value = "PropertyBrush.Nominal"
$ws.Range("B865").Value = "显示名称"
$ws.Range("C865").Value = "Nominal"
$ws.Range("A866").Value = "PropertyBrush.Tag"
$ws.Range("B866").Value = "标签"
$ws.Range("C866").Value = "Tag"
$ws.Range("A867").Value = "PropertyBrush.Status"
$ws.Range("B867").Value = "状态"
$ws.Range("C867").Value = "Status"
$ws.Range("A868").Value = "PropertyBrush.VeterancyPercentage"
$ws.Range("B868").Value = "经验等级"
$ws.Range("C868").Value = "Veterancy"
$ws.Range("A869").Value = "PropertyBrush.Group"
$ws.Range("B869").Value = "小组"
$ws.Range("C869").Value = "Group"
$ws.Range("A870").Value = "PropertyBrush.IsAboveGround"
$ws.Range("B870").Value = "在桥梁上"
$ws.Range("C870").Value = "IsAboveGround"
$ws.Range("A871").Value = "PropertyBrush.AutoNORecruitType"
$ws.Range("B871").Value = "重组A"
$ws.Range("C871").Value = "AutoNORecruitType"
$ws.Range("A872").Value = "PropertyBrush.AutoYESRecruitType"
$ws.Range("B872").Value = "重组B"
$ws.Range("C872").Value = "AutoYESRecruitType"
$ws.Range("A873").Value = "PropertyBrush.FollowsIndex"
$ws.Range("B873").Value = "跟随ID"
$ws.Range("C873").Value = "FollowsIndex"

# Update selection to reflect final cursor position
[void]$ws.Range("B860").Select()
